$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 220, shifting existing rows 220:339 down to 221:340
$ws.Rows.Item(220).Insert()

# Populate the new row 220 with the same constant values as its neighboring rows
# plus the new data-specific values from the diff.
$ws.Cells.Item(220, 1).Value = 1
$ws.Cells.Item(220, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(220, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(220, 4).Value = 44839
$ws.Cells.Item(220, 5).Value = 15
$ws.Cells.Item(220, 6).Value = 100114013
$ws.Cells.Item(220, 7).Value = "Zanahoria"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 70
$ws.Cells.Item(220, 11).Value = 21000
$ws.Cells.Item(220, 12).Value = 22000
$ws.Cells.Item(220, 13).Value = 21500
$ws.Cells.Item(220, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(220, 15).Value = "Valle de Camiña"
$ws.Cells.Item(220, 16).Value = 860
$ws.Cells.Item(220, 17).Value = 25
$ws.Cells.Item(220, 18).Value = "Hortaliza"
